$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), pushing the existing
# N/O/P ("Late" / heading / "Outstanding") columns one to the right.
$ws.Columns("N:N").Insert() | Out-Null

# Match the inserted column's width to the existing "11-wide" column (M),
# just without the bestFit flag - mirrors the target width="11" custom col.
$ws.Columns("N:N").ColumnWidth = 10.17

# Make "Repayment schedule" the active sheet/tab with cell J14 selected.
$ws.Activate() | Out-Null
$ws.Range("J14").Select() | Out-Null
